# Update the "Out of PO" roster table (A2:C19) with the newly uploaded data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Dyson Daniels", "PG,SG", "Atlanta Hawks"),
    @("Tari Eason", "SF,PF", "Houston Rockets"),
    @("De'Andre Hunter", "SF,PF", "Atlanta Hawks"),
    @("Josh Hart", "SF,PF", "New York Knicks"),
    @("Deandre Ayton", "C", "Portland Trail Blazers"),
    @("Myles Turner", "C", "Indiana Pacers"),
    @("Domantas Sabonis", "C", "Sacramento Kings"),
    @("Victor Wembanyama", "C", "San Antonio Spurs"),
    @("Bradley Beal", "PG,SG,SF", "Phoenix Suns"),
    @("Michael Porter Jr.", "SF,PF", "Denver Nuggets"),
    @("Payton Pritchard", "PG", "Boston Celtics"),
    @("Donovan Mitchell", "PG,SG", "Cleveland Cavaliers"),
    @("Malik Beasley", "SG", "Detroit Pistons"),
    @("Kristaps Porzingis", "PF,C", "Boston Celtics"),
    @("P.J. Washington", "PF", "Dallas Mavericks"),
    @("Robert Williams III", "C", "Portland Trail Blazers"),
    @("Cam Thomas", "SG,SF", "Brooklyn Nets"),
    @("Jamal Murray", "PG,SG", "Denver Nuggets")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}
